$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Todas las bases homologadas para unirse" - standardize the municipality
# column so it carries the official INEGI municipality key (CVE_MUN) instead
# of the municipality's display name, so this table can be joined to the
# other homologated bases on that key.
#
# Map of old column-A text (header + municipality names) -> new CVE_MUN text.
$muniMap = @{
    "Municipio" = "CVE_MUN"
    "Acatlán" = "13001"
    "Acaxochitlán" = "13002"
    "Actopan" = "13003"
    "Agua Blanca de Iturbide" = "13004"
    "Ajacuba" = "13005"
    "Alfajayucan" = "13006"
    "Almoloya" = "13007"
    "Apan" = "13008"
    "Atitalaquia" = "13010"
    "Atlapexco" = "13011"
    "Atotonilco de Tula" = "13013"
    "Atotonilco el Grande" = "13012"
    "Calnali" = "13014"
    "Cardonal" = "13015"
    "Chapantongo" = "13017"
    "Chapulhuacán" = "13018"
    "Chilcuautla" = "13019"
    "Cuautepec de Hinojosa" = "13016"
    "El Arenal" = "13009"
    "Eloxochitlán" = "13020"
    "Emiliano Zapata" = "13021"
    "Epazoyucan" = "13022"
    "Francisco I. Madero" = "13023"
    "Huasca de Ocampo" = "13024"
    "Huautla" = "13025"
    "Huazalingo" = "13026"
    "Huehuetla" = "13027"
    "Huejutla de Reyes" = "13028"
    "Huichapan" = "13029"
    "Ixmiquilpan" = "13030"
    "Jacala de Ledezma" = "13031"
    "Jaltocán" = "13032"
    "Juárez Hidalgo" = "13033"
    "La Misión" = "13040"
    "Lolotla" = "13034"
    "Metepec" = "13035"
    "Metztitlán" = "13037"
    "Mineral de la Reforma" = "13051"
    "Mineral del Chico" = "13038"
    "Mineral del Monte" = "13039"
    "Mixquiahuala de Juárez" = "13041"
    "Molango de Escamilla" = "13042"
    "Nicolás Flores" = "13043"
    "Nopala de Villagrán" = "13044"
    "Omitlán de Juárez" = "13045"
    "Pachuca de Soto" = "13048"
    "Pacula" = "13047"
    "Pisaflores" = "13049"
    "Progreso de Obregón" = "13050"
    "San Agustín Metzquititlán" = "13036"
    "San Agustín Tlaxiaca" = "13052"
    "San Bartolo Tutotepec" = "13053"
    "San Felipe Orizatlán" = "13046"
    "San Salvador" = "13054"
    "Santiago de Anaya" = "13055"
    "Santiago Tulantepec de Lugo Guerrero" = "13056"
    "Singuilucan" = "13057"
    "Tasquillo" = "13058"
    "Tecozautla" = "13059"
    "Tenango de Doria" = "13060"
    "Tepeapulco" = "13061"
    "Tepehuacán de Guerrero" = "13062"
    "Tepeji del Río de Ocampo" = "13063"
    "Tepetitlán" = "13064"
    "Tetepango" = "13065"
    "Tezontepec de Aldama" = "13067"
    "Tianguistengo" = "13068"
    "Tizayuca" = "13069"
    "Tlahuelilpan" = "13070"
    "Tlahuiltepa" = "13071"
    "Tlanalapa" = "13072"
    "Tlanchinol" = "13073"
    "Tlaxcoapan" = "13074"
    "Tolcayuca" = "13075"
    "Tula de Allende" = "13076"
    "Tulancingo de Bravo" = "13077"
    "Villa de Tezontepec" = "13066"
    "Xochiatipan" = "13078"
    "Xochicoatlán" = "13079"
    "Yahualica" = "13080"
    "Zacualtipán de Ángeles" = "13081"
    "Zapotlán de Juárez" = "13082"
    "Zempoala" = "13083"
    "Zimapán" = "13084"
}

# Column A holds the header in row 1 and one municipality per row below it.
# Force the column to Text format first so the numeric-looking keys (e.g.
# "13001") get written and stored as text, matching the source data, rather
# than being auto-converted to numbers.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$colA = $ws.Range("A1:A" + $lastRow)
$colA.NumberFormat = "@"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    # .Value2 reads back reliably in this runtime; .Value is used for the
    # write so the text format above is honoured.
    $current = [string]$cell.Value2
    if ($muniMap.ContainsKey($current)) {
        $cell.Value = $muniMap[$current]
    }
}
